# Updated cryptos list (prices + 1h volume change) from the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.364.87'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.311.25'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.86'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.75'
$ws.Range("E6").Value = '  -5.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.586'
$ws.Range("E7").Value = '  -3.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.305.11'
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("E10").Value = '  -2.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.95'
$ws.Range("E11").Value = '  -12.29%  '
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.94'
$ws.Range("E14").Value = '  -2.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.831.62'
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.307.43'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.116'
$ws.Range("E17").Value = '  -1.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.154.53'
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.39'
$ws.Range("E19").Value = '  -2.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.14'
$ws.Range("E20").Value = '  -0.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.953'
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '378.50'
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.16'
$ws.Range("E23").Value = '  +7.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.26'
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.11'
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.70'
$ws.Range("E26").Value = '  -4.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.16'
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  -4.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.13'
$ws.Range("E30").Value = '  -4.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.76'
$ws.Range("E31").Value = '  -1.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '628.28'
$ws.Range("E32").Value = '  -4.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.55'
$ws.Range("E33").Value = '  -3.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.14'
$ws.Range("E34").Value = '  -2.03%  '
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.12'
$ws.Range("E36").Value = '  -4.21%  '
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("E38").Value = '  -2.39%  '
$ws.Range("E39").Value = '  -4.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0754'
$ws.Range("E40").Value = '  +5.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("E41").Value = '  -0.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.20'
$ws.Range("E42").Value = '  +9.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.125'
$ws.Range("E43").Value = '  -1.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.59'
$ws.Range("E44").Value = '  +3.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.886.88'
$ws.Range("E45").Value = '  -1.00%  '
$ws.Range("E46").Value = '  +0.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0397'
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.61'
$ws.Range("E48").Value = '  -4.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.02'
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.72'
$ws.Range("E50").Value = '  +1.61%  '
$ws.Range("E51").Value = '  -2.04%  '
